# Merge the two trailing runs of the "Main running" paragraph on slide 6
# (TextBox 2) into a single run, matching the commit's simplification of
# the description text ("...it needs are." no longer a separate run).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$fullText = $tr.Text
$boldLead = "Main running"

$tailStart = $boldLead.Length + 1
$tailLength = $fullText.Length - $boldLead.Length
$tail = $tr.Characters($tailStart, $tailLength)
$tailTargetText = $tail.Text

# Force the run boundary to be rebuilt: briefly set the tail to a
# placeholder value, then restore the desired merged text so the engine
# collapses the previously separate runs ("... it " / "needs are.") into
# a single run instead of treating the identical final text as a no-op.
$tail.Text = "__tmp__"
$tail2 = $tr.Characters($tailStart, 7)
$tail2.Text = $tailTargetText
